$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1345.3636
$ws.Range("I6").Value = 479.9
$ws.Range("J6").Value = 10000
$ws.Range("K6").Value = 1439.7
$ws.Range("L6").Value = 30000
$ws.Range("M6").Value = -1327.7
$ws.Range("N6").Value = -30224

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 25799.8
$ws.Range("I8").Value = 19499.5
$ws.Range("K8").Value = 58498.5
$ws.Range("M8").Value = -58359.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 16557.666
$ws.Range("I20").Value = 16557.666
$ws.Range("K20").Value = 16557.666
$ws.Range("M20").Value = -16327.666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H35").Value = 16557.666
$ws.Range("I35").Value = 16557.666
$ws.Range("K35").Value = 16557.666
$ws.Range("M35").Value = -16178.666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1593.7778
$ws.Range("I86").Value = 1716.6666
$ws.Range("J86").Value = 1348
$ws.Range("K86").Value = 1716.6666
$ws.Range("L86").Value = 1348
$ws.Range("M86").Value = -593.6666
$ws.Range("N86").Value = -3594

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 1593.7778
$ws.Range("I89").Value = 1716.6666
$ws.Range("J89").Value = 1348
$ws.Range("K89").Value = 8583.333
$ws.Range("L89").Value = 6740
$ws.Range("M89").Value = -2967.333000000001
$ws.Range("N89").Value = -17972

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 4434.0645
$ws.Range("I98").Value = 2778.4
$ws.Range("K98").Value = 2778.4
$ws.Range("M98").Value = -1280.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 4434.0645
$ws.Range("I122").Value = 2778.4
$ws.Range("K122").Value = 8335.2
$ws.Range("M122").Value = -5885.200000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3100.1584
$ws.Range("J138").Value = 3522.1692
$ws.Range("L138").Value = 10566.5076
$ws.Range("N138").Value = -20846.5076

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2455.3
$ws.Range("I61").Value = 2087.875
$ws.Range("J61").Value = 3925
$ws.Range("K61").Value = 2087.875
$ws.Range("L61").Value = 3925
$ws.Range("M61").Value = -1875.875
$ws.Range("N61").Value = -4349

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H124").Value = 23429
$ws.Range("J124").Value = 23429
$ws.Range("L124").Value = 23429
$ws.Range("N124").Value = -33249

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1797.8704
$ws.Range("I132").Value = 1141.975
$ws.Range("J132").Value = 3671.8572
$ws.Range("K132").Value = 3425.925
$ws.Range("L132").Value = 11015.5716
$ws.Range("M132").Value = -895.9249999999997
$ws.Range("N132").Value = -16075.5716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2455.3
$ws.Range("I136").Value = 2087.875
$ws.Range("J136").Value = 3925
$ws.Range("K136").Value = 6263.625
$ws.Range("L136").Value = 11775
$ws.Range("M136").Value = -3713.625
$ws.Range("N136").Value = -16875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 6881.25
$ws.Range("I22").Value = 7850
$ws.Range("J22").Value = 100
$ws.Range("K22").Value = 7850
$ws.Range("L22").Value = 100
$ws.Range("M22").Value = -7677
$ws.Range("N22").Value = -446

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2171.7693
$ws.Range("I107").Value = 2384.818
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 2384.818
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = -464.8180000000002
$ws.Range("N107").Value = -4840

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3557.1428
$ws.Range("I134").Value = 3580
$ws.Range("K134").Value = 10740
$ws.Range("M134").Value = -8205

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 10000
$ws.Range("J4").Value = 10000
$ws.Range("L4").Value = 10000
$ws.Range("N4").Value = -10224

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 6382.2
$ws.Range("I10").Value = 2977.75
$ws.Range("K10").Value = 2977.75
$ws.Range("M10").Value = -2838.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 34500
$ws.Range("I23").Value = 34500
$ws.Range("K23").Value = 34500
$ws.Range("M23").Value = -34260

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H27").Value = 34500
$ws.Range("I27").Value = 34500
$ws.Range("K27").Value = 34500
$ws.Range("M27").Value = -34308

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1252.3182
$ws.Range("I58").Value = 1323.2106
$ws.Range("J58").Value = 803.3333
$ws.Range("K58").Value = 1323.2106
$ws.Range("L58").Value = 803.3333
$ws.Range("M58").Value = -1120.2106
$ws.Range("N58").Value = -1209.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 34980
$ws.Range("J60").Value = 37675
$ws.Range("L60").Value = 37675
$ws.Range("N60").Value = -38697

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2245
$ws.Range("I122").Value = 2000.6666
$ws.Range("J122").Value = 2831.4
$ws.Range("K122").Value = 6001.9998
$ws.Range("L122").Value = 8494.2
$ws.Range("M122").Value = -3551.9998
$ws.Range("N122").Value = -13394.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1252.3182
$ws.Range("I136").Value = 1323.2106
$ws.Range("J136").Value = 803.3333
$ws.Range("K136").Value = 3969.6318
$ws.Range("L136").Value = 2409.9999
$ws.Range("M136").Value = -1419.6318
$ws.Range("N136").Value = -7509.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 526958.75
$ws.Range("I113").Value = 1111725.2
$ws.Range("J113").Value = 668.9
$ws.Range("K113").Value = 3335175.6
$ws.Range("L113").Value = 2006.7
$ws.Range("M113").Value = -3333005.6
$ws.Range("N113").Value = -6346.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2804
$ws.Range("I140").Value = 1610
$ws.Range("J140").Value = 4215.091
$ws.Range("K140").Value = 4830
$ws.Range("L140").Value = 12645.273
$ws.Range("M140").Value = 350
$ws.Range("N140").Value = -23005.273

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 8050.75
$ws.Range("I21").Value = 4400
$ws.Range("J21").Value = 9267.667
$ws.Range("K21").Value = 4400
$ws.Range("L21").Value = 9267.667
$ws.Range("M21").Value = -4227
$ws.Range("N21").Value = -9613.667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H30").Value = 8050.75
$ws.Range("I30").Value = 4400
$ws.Range("J30").Value = 9267.667
$ws.Range("K30").Value = 4400
$ws.Range("L30").Value = 9267.667
$ws.Range("M30").Value = -4295
$ws.Range("N30").Value = -9477.667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2519.0527
$ws.Range("I102").Value = 2162.4
$ws.Range("J102").Value = 2915.3333
$ws.Range("K102").Value = 2162.4
$ws.Range("L102").Value = 2915.3333
$ws.Range("M102").Value = -540.4000000000001
$ws.Range("N102").Value = -6159.3333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2338.875
$ws.Range("I126").Value = 1452.75
$ws.Range("J126").Value = 3225
$ws.Range("K126").Value = 4358.25
$ws.Range("L126").Value = 9675
$ws.Range("M126").Value = -1888.25
$ws.Range("N126").Value = -14615

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2086.5454
$ws.Range("I132").Value = 1749.5555
$ws.Range("J132").Value = 3603
$ws.Range("K132").Value = 5248.666499999999
$ws.Range("L132").Value = 10809
$ws.Range("M132").Value = -2718.666499999999
$ws.Range("N132").Value = -15869

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 990.9091
$ws.Range("I22").Value = 457.14285
$ws.Range("J22").Value = 1925
$ws.Range("K22").Value = 457.14285
$ws.Range("L22").Value = 1925
$ws.Range("M22").Value = -162.14285
$ws.Range("N22").Value = -2515

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 6502.3335
$ws.Range("J23").Value = 9003.5
$ws.Range("L23").Value = 9003.5
$ws.Range("N23").Value = -9463.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 990.9091
$ws.Range("I27").Value = 457.14285
$ws.Range("J27").Value = 1925
$ws.Range("K27").Value = 457.14285
$ws.Range("L27").Value = 1925
$ws.Range("M27").Value = -350.14285
$ws.Range("N27").Value = -2139

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2987.5
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 3152.0833
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 3152.0833
$ws.Range("M40").Value = -1864
$ws.Range("N40").Value = -3424.0833

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 1000.6667
$ws.Range("I2").Value = 1002
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 1002
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = -890
$ws.Range("N2").Value = -1224

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 23333.666
$ws.Range("I11").Value = 30000
$ws.Range("K11").Value = 30000
$ws.Range("M11").Value = -29858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 10419018
$ws.Range("I122").Value = 11365700
$ws.Range("K122").Value = 34097100
$ws.Range("M122").Value = -34094650
